$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "end" time for the row 50 entry (E50) from 12:00 to 13:00
$ws.Range("E50").Value = 0.54166666666666663

# Update the selected cell to E51 (matches the sheetView selection change)
$ws.Range("E51").Select()
